$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Имя заявки" (${experiment.firstName}) column is column E on this
# sheet. Delete it entirely - Excel shifts every following column one to
# the left, so the former last column M becomes column L, and every row's
# spans/merges shrink from 1:13 to 1:12 automatically.
$ws.Range("E1").EntireColumn.Delete()

# The A1/A3/A4/A7 cell comments hold JasperXLS template directives that
# hard-code the sheet's last column letter (lastCell="M..."). Now that the
# last column is L instead of M, update those references.
foreach ($addr in @("A1", "A3", "A4", "A7")) {
    $cmt = $ws.Range($addr).Comment
    $newText = $cmt.Text() -replace 'lastCell="M', 'lastCell="L'
    [void]$cmt.Text($newText)
}

# Match the author's final selection in the saved file.
[void]$ws.Range("A8:L8").Select()
